$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("C2").Value = 65
$ws.Range("G2").Value = 80
$ws.Range("J2").Value = 108
$ws.Range("E3").Value = 129
$ws.Range("C6").Value = 439
$ws.Range("D6").Value = 372
$ws.Range("E6").Value = 413
$ws.Range("F6").Value = 463
$ws.Range("G6").Value = 413
$ws.Range("H6").Value = 412
$ws.Range("I6").Value = 462
$ws.Range("C7").Value = 586
$ws.Range("D7").Value = 582
$ws.Range("E7").Value = 617
$ws.Range("F7").Value = 667
$ws.Range("G7").Value = 624
$ws.Range("H7").Value = 656
$ws.Range("I7").Value = 775
$ws.Range("J7").Value = 708

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("E3").Value = 6
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 11

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("E5").Value = 16
$ws.Range("I5").Value = 9
$ws.Range("E6").Value = 18
$ws.Range("I6").Value = 17

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("C6").Value = 39
$ws.Range("F6").Value = 37
$ws.Range("C7").Value = 44
$ws.Range("F7").Value = 51

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("E5").Value = 11
$ws.Range("D7").Value = 9
$ws.Range("E8").Value = 44
$ws.Range("G8").Value = 30
$ws.Range("G27").Value = 9
$ws.Range("C28").Value = 44
$ws.Range("F28").Value = 51
$ws.Range("J36").Value = 41
$ws.Range("E50").Value = 18
$ws.Range("I50").Value = 17
$ws.Range("C53").Value = 51
$ws.Range("C56").Value = 4
$ws.Range("G65").Value = 18
$ws.Range("H65").Value = 16
$ws.Range("I70").Value = 18
$ws.Range("C77").Value = 23
$ws.Range("E78").Value = 7
$ws.Range("F79").Value = 8
$ws.Range("E92").Value = 4
$ws.Range("C98").Value = 586
$ws.Range("D98").Value = 582
$ws.Range("E98").Value = 617
$ws.Range("F98").Value = 667
$ws.Range("G98").Value = 624
$ws.Range("H98").Value = 656
$ws.Range("I98").Value = 775
$ws.Range("J98").Value = 708

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("E4").Value = 4
$ws.Range("E5").Value = 7

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("C6").Value = 33
$ws.Range("C7").Value = 51

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 4

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 8

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 8
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = 16

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("C6").Value = 15
$ws.Range("C7").Value = 23

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("G2").Value = 2
$ws.Range("G5").Value = 9

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 9

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("C2").Value = 2
$ws.Range("C5").Value = 4

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I5").Value = 16
$ws.Range("I6").Value = 18

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E5").Value = 34
$ws.Range("G5").Value = 22
$ws.Range("E6").Value = 44
$ws.Range("G6").Value = 30
